$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# New header cell for the SnapshotURL column (AZ), mirroring the existing
# header row's last cell style (fill/font used on AT1:AY1).
$ws.Range("AZ1").Value = "SnapshotURL"
$ws.Range("AY1").Copy()
$ws.Range("AZ1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the new column for every data row (2-215) with the snapshot URL.
$url = "https://education-profiles.org/central-and-southern-asia/afghanistan/~inclusion"
$ws.Range("AZ2:AZ215").Value = $url

# Match the saved selection state: AZ2:AZ215 selected with AZ2 active.
$ws.Range("AZ2:AZ215").Select()
